$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32, pushing existing rows 32:71 down to 33:72
$ws.Rows("32:32").Insert()

# Fill in the new row 32 with the new data
$ws.Range("A32").Value = 4
$ws.Range("B32").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C32").Value = "Los Lagos"
$ws.Range("D32").Value = 44540
$ws.Range("E32").Value = 10
$ws.Range("F32").Value = 100112026
$ws.Range("G32").Value = "Haba"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 120
$ws.Range("K32").Value = 12000
$ws.Range("L32").Value = 12000
$ws.Range("M32").Value = 12000
$ws.Range("N32").Value = "$/saco 25 kilos"
$ws.Range("O32").Value = "Región de La Araucanía"
$ws.Range("P32").Value = 480
$ws.Range("Q32").Value = 25
$ws.Range("R32").Value = "Hortaliza"
